$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 158.3
$ws.Range("I33").Value = 158.3
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 158.3
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 70.69999999999999
$ws.Range("N33").ClearContents()

$ws.Range("H40").Value = 2571.1428
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1825

$ws.Range("H64").Value = 2398.4
$ws.Range("I64").Value = 2399
$ws.Range("J64").Value = 2397.5
$ws.Range("K64").Value = 2399
$ws.Range("L64").Value = 2397.5
$ws.Range("M64").Value = -2151
$ws.Range("N64").Value = -2893.5

$ws.Range("H67").Value = 2398.4
$ws.Range("I67").Value = 2399
$ws.Range("J67").Value = 2397.5
$ws.Range("K67").Value = 2399
$ws.Range("L67").Value = 2397.5
$ws.Range("M67").Value = -1541
$ws.Range("N67").Value = -4113.5

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H92").Value = 328.70587
$ws.Range("I92").Value = 212.5
$ws.Range("K92").Value = 212.5
$ws.Range("M92").Value = 1035.5

$ws.Range("H98").Value = 1586.125
$ws.Range("I98").Value = 1637.2307
$ws.Range("J98").Value = 1364.6666
$ws.Range("K98").Value = 1637.2307
$ws.Range("L98").Value = 1364.6666
$ws.Range("M98").Value = -139.2307000000001
$ws.Range("N98").Value = -4360.6666

$ws.Range("H122").Value = 1586.125
$ws.Range("I122").Value = 1637.2307
$ws.Range("J122").Value = 1364.6666
$ws.Range("K122").Value = 4911.6921
$ws.Range("L122").Value = 4093.9998
$ws.Range("M122").Value = -2461.6921
$ws.Range("N122").Value = -8993.9998

$ws.Range("H132").Value = 2845
$ws.Range("I132").Value = 1533.2354
$ws.Range("K132").Value = 4599.706200000001
$ws.Range("M132").Value = -2069.706200000001

$ws.Range("H135").Value = 1692.3846
$ws.Range("I135").Value = 1270.6666
$ws.Range("K135").Value = 11435.9994
$ws.Range("M135").Value = -8900.999400000001

$ws.Range("H137").Value = 1349.2963
$ws.Range("I137").Value = 899.65
$ws.Range("K137").Value = 2698.95
$ws.Range("M137").Value = -148.9499999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2037.6316
$ws.Range("I45").Value = 1486.7858
$ws.Range("K45").Value = 1486.7858
$ws.Range("M45").Value = -1109.7858

$ws.Range("H61").Value = 9000.714
$ws.Range("I61").Value = 9000.714
$ws.Range("K61").Value = 9000.714
$ws.Range("M61").Value = -8788.714

$ws.Range("H74").Value = 2031.5
$ws.Range("I74").Value = 1245.2
$ws.Range("K74").Value = 1245.2
$ws.Range("M74").Value = -371.2

$ws.Range("H77").Value = 2031.5
$ws.Range("I77").Value = 1245.2
$ws.Range("K77").Value = 6226
$ws.Range("M77").Value = -1858

$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 5000
$ws.Range("N101").Value = -11490

$ws.Range("H122").Value = 12720.77
$ws.Range("I122").Value = 9423.478999999999
$ws.Range("K122").Value = 28270.437
$ws.Range("M122").Value = -25820.437

$ws.Range("H132").Value = 8173.875
$ws.Range("I132").Value = 11899
$ws.Range("K132").Value = 35697
$ws.Range("M132").Value = -33167

$ws.Range("H136").Value = 9000.714
$ws.Range("I136").Value = 9000.714
$ws.Range("K136").Value = 27002.142
$ws.Range("M136").Value = -24452.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 8550
$ws.Range("J15").Value = 8550
$ws.Range("L15").Value = 8550
$ws.Range("N15").Value = -9004

$ws.Range("H99").Value = 2104.4285
$ws.Range("I99").Value = 2994
$ws.Range("J99").Value = 1437.25
$ws.Range("K99").Value = 2994
$ws.Range("L99").Value = 1437.25
$ws.Range("M99").Value = -1496
$ws.Range("N99").Value = -4433.25

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws.Range("H130").Value = 57655.555
$ws.Range("J130").Value = 57655.555
$ws.Range("L130").Value = 57655.555
$ws.Range("N130").Value = -67695.55499999999

$ws.Range("H134").Value = 3377.6365
$ws.Range("I134").Value = 2628.2222
$ws.Range("K134").Value = 7884.6666
$ws.Range("M134").Value = -5349.6666

$ws.Range("H140").Value = 80444.5
$ws.Range("I140").Value = 89999
$ws.Range("J140").Value = 70890
$ws.Range("K140").Value = 89999
$ws.Range("L140").Value = 70890
$ws.Range("M140").Value = -84819
$ws.Range("N140").Value = -81250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 103499.5
$ws.Range("I62").Value = 4666.3335
$ws.Range("K62").Value = 4666.3335
$ws.Range("M62").Value = -4042.3335

$ws.Range("H65").Value = 103499.5
$ws.Range("I65").Value = 4666.3335
$ws.Range("K65").Value = 23331.6675
$ws.Range("M65").Value = -20211.6675

$ws.Range("H105").Value = 1027.5
$ws.Range("I105").Value = 1027.5
$ws.Range("K105").Value = 1027.5
$ws.Range("M105").Value = 719.5

$ws.Range("H134").Value = 3196.6875
$ws.Range("I134").Value = 2832.6365
$ws.Range("K134").Value = 8497.9095
$ws.Range("M134").Value = -5962.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2872163.8
$ws.Range("I4").Value = 3648292
$ws.Range("J4").Value = 488.7
$ws.Range("K4").Value = 10944876
$ws.Range("L4").Value = 1466.1
$ws.Range("M4").Value = -10944764
$ws.Range("N4").Value = -1690.1

$ws.Range("H29").Value = 521.8
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 627.25
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 1881.75
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -2435.75

$ws.Range("H114").Value = 1263.4
$ws.Range("I114").Value = 1059
$ws.Range("J114").Value = 1399.6666
$ws.Range("K114").Value = 3177
$ws.Range("L114").Value = 4198.9998
$ws.Range("M114").Value = 77
$ws.Range("N114").Value = -10706.9998

$ws.Range("H117").Value = 1145.875
$ws.Range("I117").Value = 650.3333
$ws.Range("J117").Value = 1443.2
$ws.Range("K117").Value = 1950.9999
$ws.Range("L117").Value = 4329.6
$ws.Range("M117").Value = 1491.0001
$ws.Range("N117").Value = -11213.6

$ws.Range("H140").Value = 1308.8
$ws.Range("I140").Value = 1045.2858
$ws.Range("K140").Value = 3135.8574
$ws.Range("M140").Value = 2044.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4996.2
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 4996.2
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984

$ws.Range("H100").Value = 39999
$ws.Range("J100").Value = 39999
$ws.Range("L100").Value = 39999
$ws.Range("N100").Value = -42163

$ws.Range("H101").Value = 44999.668
$ws.Range("J101").Value = 44999.668
$ws.Range("L101").Value = 44999.668
$ws.Range("N101").Value = -51489.668

$ws.Range("H122").Value = 74639.21000000001
$ws.Range("I122").Value = 2618.1
$ws.Range("J122").Value = 254692
$ws.Range("K122").Value = 7854.299999999999
$ws.Range("L122").Value = 764076
$ws.Range("M122").Value = -5404.299999999999
$ws.Range("N122").Value = -768976

$ws.Range("H132").Value = 3123.0908
$ws.Range("I132").Value = 3132
$ws.Range("K132").Value = 9396
$ws.Range("M132").Value = -6866

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1595.5
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H40").Value = 3536.68
$ws.Range("I40").Value = 3551.5833
$ws.Range("K40").Value = 3551.5833
$ws.Range("M40").Value = -3415.5833

$ws.Range("H122").Value = 5794.2
$ws.Range("I122").Value = 4419.4287
$ws.Range("J122").Value = 6997.125
$ws.Range("K122").Value = 13258.2861
$ws.Range("L122").Value = 20991.375
$ws.Range("M122").Value = -10808.2861
$ws.Range("N122").Value = -25891.375

$ws.Range("H132").Value = 54146.316
$ws.Range("I132").Value = 56960
$ws.Range("K132").Value = 170880
$ws.Range("M132").Value = -168350

$ws.Range("H140").Value = 100214
$ws.Range("J140").Value = 100214
$ws.Range("L140").Value = 100214
$ws.Range("N140").Value = -110574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H93").Value = 45999.332
$ws.Range("J93").Value = 45999.332
$ws.Range("L93").Value = 45999.332
$ws.Range("N93").Value = -50991.332

$ws.Range("H103").Value = 17824.4
$ws.Range("J103").Value = 17824.4
$ws.Range("L103").Value = 17824.4
$ws.Range("N103").Value = -20168.4

$ws.Range("H122").Value = 2840.8
$ws.Range("I122").Value = 2531.2
$ws.Range("K122").Value = 7593.599999999999
$ws.Range("M122").Value = -5143.599999999999
